# Auto-generated Excel COM-interop edit script
# Applies the scheduled-runner value refresh described in the commit diff
# across sheets ALC, ARM, BSM, CRP, CUL, GSM, LTW (WVR untouched).

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# ALC!row 4 - Root Rush
$ws.Range("H4").Value = 125.25
$ws.Range("I4").Value = 100.5
$ws.Range("J4").Value = 150
$ws.Range("K4").Value = 100.5
$ws.Range("L4").Value = 150
$ws.Range("M4").Value = 13.5
$ws.Range("N4").Value = -378

# ALC!row 8 - On the Drip
$ws.Range("H8").Value = 77114.234
$ws.Range("I8").Value = 77114.234
$ws.Range("K8").Value = 231342.702
$ws.Range("M8").Value = -231203.702

# ALC!row 62 - The Mustache Suits Him
$ws.Range("H62").Value = 65013.41
$ws.Range("I62").Value = 76237.78999999999
$ws.Range("J62").Value = 12633
$ws.Range("K62").Value = 76237.78999999999
$ws.Range("L62").Value = 12633
$ws.Range("M62").Value = -75613.78999999999
$ws.Range("N62").Value = -13881

# ALC!row 65 - Forgery of Convenience (L)
$ws.Range("H65").Value = 65013.41
$ws.Range("I65").Value = 76237.78999999999
$ws.Range("J65").Value = 12633
$ws.Range("K65").Value = 381188.95
$ws.Range("L65").Value = 63165
$ws.Range("M65").Value = -378068.95
$ws.Range("N65").Value = -69405

$ws = $wb.Worksheets.Item("ARM")
# ARM!row 10 - Bronzed and Burnt
$ws.Range("H10").Value = 0
$ws.Range("I10").Value = 0
$ws.Range("J10").Value = 0
$ws.Range("K10").Value = 0
$ws.Range("L10").ClearContents()
$ws.Range("M10").ClearContents()
$ws.Range("N10").Value = 0

# ARM!row 32 - Ingot We Trust
$ws.Range("H32").Value = 3960.66
$ws.Range("I32").Value = 3986.5251
$ws.Range("J32").Value = 1400
$ws.Range("K32").Value = 3986.5251
$ws.Range("L32").Value = 1400
$ws.Range("M32").Value = -3699.5251
$ws.Range("N32").Value = -1974

# ARM!row 37 - Get Shirty
$ws.Range("H37").Value = 8203.5
$ws.Range("J37").Value = 8203.5
$ws.Range("L37").Value = 8203.5
$ws.Range("N37").Value = -8749.5

# ARM!row 44 - Very Slow Array
$ws.Range("H44").Value = 22049
$ws.Range("J44").Value = 22049
$ws.Range("L44").Value = 22049
$ws.Range("N44").Value = -23025

# ARM!row 55 - Employee Retention
$ws.Range("H55").Value = 25053
$ws.Range("J55").Value = 25053
$ws.Range("L55").Value = 25053
$ws.Range("N55").Value = -25683

$ws = $wb.Worksheets.Item("BSM")
# BSM!row 107 - The Gold Experience
$ws.Range("H107").Value = 10982.23
$ws.Range("I107").Value = 830.75
$ws.Range("J107").Value = 132800
$ws.Range("K107").Value = 830.75
$ws.Range("L107").Value = 132800
$ws.Range("M107").Value = 1089.25
$ws.Range("N107").Value = -136640

$ws = $wb.Worksheets.Item("CRP")
# CRP!row 22 - Driving Up the Wall
$ws.Range("H22").Value = 1123.6666
$ws.Range("I22").Value = 1385.25
$ws.Range("J22").Value = 600.5
$ws.Range("K22").Value = 1385.25
$ws.Range("L22").Value = 600.5
$ws.Range("M22").Value = -1035.25
$ws.Range("N22").Value = -1300.5

# CRP!row 31 - Wall Not Found
$ws.Range("H31").Value = 103884.3
$ws.Range("I31").Value = 0
$ws.Range("J31").Value = 103884.3
$ws.Range("K31").Value = 0
$ws.Range("L31").ClearContents()
$ws.Range("M31").Value = 103884.3
$ws.Range("N31").Value = -104474.3

# CRP!row 34 - Armoires of the Rich and Famous
$ws.Range("H34").Value = 103884.3
$ws.Range("I34").Value = 0
$ws.Range("J34").Value = 103884.3
$ws.Range("K34").Value = 0
$ws.Range("L34").ClearContents()
$ws.Range("M34").Value = 103884.3
$ws.Range("N34").Value = -104288.3

$ws = $wb.Worksheets.Item("CUL")
# CUL!row 5 - What a Sap
$ws.Range("H5").Value = 2147.8333
$ws.Range("I5").Value = 1260.3125
$ws.Range("J5").Value = 2591.5938
$ws.Range("K5").Value = 3780.9375
$ws.Range("L5").Value = 7774.7814
$ws.Range("M5").Value = -3668.9375
$ws.Range("N5").Value = -7998.7814

# CUL!row 13 - Fishy Revelations
$ws.Range("H13").Value = 300.5
$ws.Range("I13").Value = 401
$ws.Range("J13").Value = 200
$ws.Range("K13").Value = 1203
$ws.Range("L13").Value = 600
$ws.Range("M13").Value = -1035
$ws.Range("N13").Value = -936

# CUL!row 32 - Convalescence Precedes Essence
$ws.Range("H32").Value = 3400
$ws.Range("J32").Value = 3500
$ws.Range("L32").Value = 10500
$ws.Range("N32").Value = -11066

# CUL!row 113 - Can't Eat Just One
$ws.Range("H113").Value = 668.76
$ws.Range("J113").Value = 547.5238000000001
$ws.Range("L113").Value = 1642.5714
$ws.Range("N113").Value = -5982.571400000001

# CUL!row 126 - Imperial Palate
$ws.Range("H126").Value = 35546.875
$ws.Range("I126").Value = 127155
$ws.Range("J126").Value = 5010.8335
$ws.Range("K126").Value = 381465
$ws.Range("L126").Value = 15032.5005
$ws.Range("M126").Value = -376525
$ws.Range("N126").Value = -24912.5005

# CUL!row 135 - Not-so-secret Ingredient
$ws.Range("H135").Value = 2147.8333
$ws.Range("I135").Value = 1260.3125
$ws.Range("J135").Value = 2591.5938
$ws.Range("K135").Value = 11342.8125
$ws.Range("L135").Value = 23324.3442
$ws.Range("M135").Value = -8807.8125
$ws.Range("N135").Value = -28394.3442

$ws = $wb.Worksheets.Item("GSM")
# GSM!row 17 - Point of Honor
$ws.Range("H17").Value = 50009
$ws.Range("I17").Value = 0
$ws.Range("J17").Value = 50009
$ws.Range("K17").Value = 0
$ws.Range("L17").ClearContents()
$ws.Range("M17").Value = 50009
$ws.Range("N17").Value = -50345

# GSM!row 80 - Needs More Prayerbell
$ws.Range("H80").Value = 3728.652
$ws.Range("I80").Value = 3697.842
$ws.Range("J80").Value = 3875
$ws.Range("K80").Value = 3697.842
$ws.Range("L80").Value = 3875
$ws.Range("M80").Value = -2699.842
$ws.Range("N80").Value = -5871

# GSM!row 83 - With a Noise That Reaches Heaven (L)
$ws.Range("H83").Value = 3728.652
$ws.Range("I83").Value = 3697.842
$ws.Range("J83").Value = 3875
$ws.Range("K83").Value = 18489.21
$ws.Range("L83").Value = 19375
$ws.Range("M83").Value = -13497.21
$ws.Range("N83").Value = -29359

# GSM!row 102 - Put the Metal to the Peddle
$ws.Range("H102").Value = 1800
$ws.Range("I102").Value = 1875
$ws.Range("J102").Value = 1500
$ws.Range("K102").Value = 1875
$ws.Range("L102").Value = 1500
$ws.Range("M102").Value = -253
$ws.Range("N102").Value = -4744

# GSM!row 113 - Copious Crystal Cannons
$ws.Range("H113").Value = 4651.1763
$ws.Range("I113").Value = 5329.091
$ws.Range("J113").Value = 3408.3333
$ws.Range("K113").Value = 5329.091
$ws.Range("L113").Value = 3408.3333
$ws.Range("M113").Value = -3159.091
$ws.Range("N113").Value = -7748.3333

$ws = $wb.Worksheets.Item("LTW")
# LTW!row 132 - Tenets of Tanning
$ws.Range("H132").Value = 2666.7183
$ws.Range("I132").Value = 1917.4717
$ws.Range("J132").Value = 4872.8335
$ws.Range("K132").Value = 5752.4151
$ws.Range("L132").Value = 14618.5005
$ws.Range("M132").Value = -3222.4151
$ws.Range("N132").Value = -19678.5005

# LTW!row 136 - Respect for Br'aax
$ws.Range("H136").Value = 3309.3572
$ws.Range("I136").Value = 1792.7084
$ws.Range("K136").Value = 5378.1252
$ws.Range("M136").Value = -2828.1252

